$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 2167.314610709427
$ws.Range("C2").Value = 1511.213283282402
$ws.Range("D2").Value = 1229.363949408971
$ws.Range("E2").Value = 2205.423038560061
$ws.Range("F2").Value = 2170.085420450942
$ws.Range("G2").Value = 1984.399082413055
$ws.Range("H2").Value = 2206.848774439267
$ws.Range("B3").Value = 2183.632415180711
$ws.Range("C3").Value = 1528.33010704203
$ws.Range("D3").Value = 1085.127383123723
$ws.Range("E3").Value = 2208.593229574046
$ws.Range("F3").Value = 2184.918646513586
$ws.Range("G3").Value = 1961.485732241424
$ws.Range("H3").Value = 2209.291490143483
$ws.Range("B4").Value = 2131.347061543535
$ws.Range("C4").Value = 1540.626423965915
$ws.Range("D4").Value = 1078.054668413406
$ws.Range("E4").Value = 2192.085558436731
$ws.Range("F4").Value = 2135.372170126485
$ws.Range("G4").Value = 1962.522550761016
$ws.Range("H4").Value = 2194.109350309947
$ws.Range("B5").Value = 2185.449126116547
$ws.Range("C5").Value = 1546.222606126279
$ws.Range("D5").Value = 1117.624230808347
$ws.Range("E5").Value = 2206.043190708929
$ws.Range("F5").Value = 2186.274067604017
$ws.Range("G5").Value = 1923.028305113323
$ws.Range("H5").Value = 2206.671397389214
$ws.Range("B6").Value = 2191.656695839558
$ws.Range("C6").Value = 1536.915018266464
$ws.Range("D6").Value = 382.5936118305735
$ws.Range("E6").Value = 2213.253418125554
$ws.Range("F6").Value = 2192.367416912643
$ws.Range("G6").Value = 1625.470157090052
$ws.Range("H6").Value = 2213.597321349207
$ws.Range("B7").Value = 2172.86595052552
$ws.Range("C7").Value = 1564.451171895376
$ws.Range("D7").Value = 572.1422959108772
$ws.Range("E7").Value = 2200.677181678597
$ws.Range("F7").Value = 2174.458556810208
$ws.Range("G7").Value = 1681.95337263134
$ws.Range("H7").Value = 2201.945082437705
$ws.Range("B8").Value = 2173.684592583933
$ws.Range("C8").Value = 1485.222641346011
$ws.Range("D8").Value = 656.094413528098
$ws.Range("E8").Value = 2196.697326503753
$ws.Range("F8").Value = 2174.349416912381
$ws.Range("G8").Value = 1665.712715699316
$ws.Range("H8").Value = 2197.242565447381
$ws.Range("B9").Value = 2187.12321845688
$ws.Range("C9").Value = 1548.913673198975
$ws.Range("D9").Value = 641.685215767798
$ws.Range("E9").Value = 2200.595623566658
$ws.Range("F9").Value = 2188.916812052127
$ws.Range("G9").Value = 1687.562082513365
$ws.Range("H9").Value = 2202.093252104228
$ws.Range("B10").Value = 1937.729501449272
$ws.Range("C10").Value = 1647.746076806086
$ws.Range("D10").Value = 1143.705851369911
$ws.Range("E10").Value = 2097.389136720578
$ws.Range("F10").Value = 1970.567967502973
$ws.Range("G10").Value = 1909.025596153625
$ws.Range("H10").Value = 2111.90225829121
$ws.Range("B11").Value = 1884.386919791947
$ws.Range("C11").Value = 1678.467432057428
$ws.Range("D11").Value = 967.7911373584814
$ws.Range("E11").Value = 2104.203142466426
$ws.Range("F11").Value = 1908.487280735242
$ws.Range("G11").Value = 1890.486000554413
$ws.Range("H11").Value = 2111.254853176406
$ws.Range("B12").Value = 1614.595651276729
$ws.Range("C12").Value = 1616.42836104628
$ws.Range("D12").Value = 267.2320334165278
$ws.Range("E12").Value = 1991.94078381052
$ws.Range("F12").Value = 1624.389466078262
$ws.Range("G12").Value = 1669.150632374289
$ws.Range("H12").Value = 1993.751378557799
$ws.Range("B13").Value = 1919.499697099683
$ws.Range("C13").Value = 1641.824615553334
$ws.Range("D13").Value = 891.4134588035242
$ws.Range("E13").Value = 2087.825279412683
$ws.Range("F13").Value = 1936.788313670181
$ws.Range("G13").Value = 1835.465197322763
$ws.Range("H13").Value = 2094.715307719508
